$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row formatting down for the 9 new data rows (118 / event "02") ---
# Rows 155-162 mirror the "normal" row style (no bottom border), like row 153.
$ws.Range("A153:F153").Copy()
$ws.Range("A155:F162").PasteSpecial(-4122)

# Row 163 is the last row of the group, mirrors the bordered style of row 154.
$ws.Range("A154:F154").Copy()
$ws.Range("A163:F163").PasteSpecial(-4122)

# --- Ensure the Script column (numeric-looking "02") is stored as text, not a number ---
$ws.Range("B155:B163").NumberFormat = "@"

# --- Row 155: Group 118, Script 02, Event heckling ---
$ws.Range("A155").Value = 118
$ws.Range("B155").Value = "02"
$ws.Range("C155").Value = "heckling"
$ws.Range("D155").Value = 3
$ws.Range("E155").Value = 5
$ws.Range("F155").Value = 4

# --- Row 156: Group 118, Script 02, Event snipping ---
$ws.Range("A156").Value = 118
$ws.Range("B156").Value = "02"
$ws.Range("C156").Value = "snipping"
$ws.Range("D156").Value = 0
$ws.Range("E156").Value = 10
$ws.Range("F156").Value = 8

# --- Row 157: Group 118, Script 02, Event whispering ---
$ws.Range("A157").Value = 118
$ws.Range("B157").Value = "02"
$ws.Range("C157").Value = "whispering"
$ws.Range("D157").Value = 5
$ws.Range("E157").Value = 7
$ws.Range("F157").Value = 9

# --- Row 158: Group 118, Script 02, Event head on table ---
$ws.Range("A158").Value = 118
$ws.Range("B158").Value = "02"
$ws.Range("C158").Value = "head on table"
$ws.Range("D158").Value = 6
$ws.Range("E158").Value = 6
$ws.Range("F158").Value = 5

# --- Row 159: Group 118, Script 02, Event drawing ---
$ws.Range("A159").Value = 118
$ws.Range("B159").Value = "02"
$ws.Range("C159").Value = "drawing"
$ws.Range("D159").Value = 0
$ws.Range("E159").Value = 10
$ws.Range("F159").Value = 9

# --- Row 160: Group 118, Script 02, Event drumming ---
$ws.Range("A160").Value = 118
$ws.Range("B160").Value = "02"
$ws.Range("C160").Value = "drumming"
$ws.Range("D160").Value = 9
$ws.Range("E160").Value = 9
$ws.Range("F160").Value = 6

# --- Row 161: Group 118, Script 02, Event locking at phone ---
$ws.Range("A161").Value = 118
$ws.Range("B161").Value = "02"
$ws.Range("C161").Value = "locking at phone"
$ws.Range("D161").Value = 1
$ws.Range("E161").Value = 9
$ws.Range("F161").Value = 10

# --- Row 162: Group 118, Script 02, Event clicking pen ---
$ws.Range("A162").Value = 118
$ws.Range("B162").Value = "02"
$ws.Range("C162").Value = "clicking pen"
$ws.Range("D162").Value = 9
$ws.Range("E162").Value = 9
$ws.Range("F162").Value = 8

# --- Row 163: Group 118, Script 02, Event "chatting " (trailing space variant) ---
$ws.Range("A163").Value = 118
$ws.Range("B163").Value = "02"
$ws.Range("C163").Value = "chatting "
$ws.Range("D163").Value = 10
$ws.Range("E163").Value = 4
$ws.Range("F163").Value = 5

# --- Update the view: scroll position + active selection moved while entering data ---
$ws.Application.ActiveWindow.ScrollRow = 142
[void]$ws.Range("C165").Select()
